# ============================================================
# IronSP_Install.docx -- "Iron Item Event Receiver added"
# ============================================================
$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 2) Append the new "Iron Event Receivers" section at the end of the
#    document (after the last paragraph, which reads "end").
# ------------------------------------------------------------------
$idx = $d.Paragraphs.Count

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Heading 1'
$cur.Range.InsertAfter('Iron Event Receivers')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter('You can register list item events on any list with the following code:')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter('include Microsoft::SharePoint')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter('class TestEventReceiver < SPItemEventReceiver')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter(([string][char]9) + 'def ItemUpdated(props)')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter(([string][char]9) + ([string][char]9) + 'self.EventFiringEnabled = false' + ([string][char]9) + ' ')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter(([string][char]9) + ([string][char]9) + 'item = props.ListItem')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter(([string][char]9) + ([string][char]9) + 'item["Title"] = "Iron Hive Event receiver"')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter(([string][char]9) + ([string][char]9) + 'item.SystemUpdate()')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter(([string][char]9) + ([string][char]9) + 'self.EventFiringEnabled = true')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter(([string][char]9) + 'end')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter('end')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter('$RUNTIME.RegisterDynamicType("TestEventReceiver", TestEventReceiver )')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter('list = SPSite.new("http://intranet/sites/IronSharePoint").RootWeb.Lists["Announcements"]')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter('# register will only  register events that has not been registered before' + ([string][char]11) + 'IronSharePoint::EventReceivers::IronItemEventReceiver.Register(list, SPEventReceiverType.ItemUpdated, SPEventReceiverSynchronization.Default, 1000, "TestEventReceiver")')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter('# IronSharePoint::EventReceivers::IronItemEventReceiver .GetAllRegistered(list)')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter('# IronSharePoint::EventReceivers::IronItemEventReceiver. Unregister(list, SPEventReceiverType.ItemUpdated, ' + ([string][char]8220) + 'TestEventReceiver' + ([string][char]8221) + ')')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'
$cur.Range.InsertAfter('# IronSharePoint::EventReceivers::IronItemEventReceiver. IsRegistered(list, SPEventReceiverType.ItemUpdated, ' + ([string][char]8220) + 'TestEventReceiver' + ([string][char]8221) + ')')

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'

$cur = $d.Paragraphs($idx)
$cur.Range.InsertParagraphAfter()
$idx++
$cur = $d.Paragraphs($idx)
$cur.Style = 'Normal'

